$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1. Update the F column (time_taken) timestamps on the "data" sheet ---
$timeValues = @(
    "2021-10-05 14:34:06.783186",
    "2021-10-05 14:34:06.783194",
    "2021-10-05 14:34:06.783207",
    "2021-10-05 14:34:06.783210",
    "2021-10-05 14:34:06.783213",
    "2021-10-05 14:34:06.783216",
    "2021-10-05 14:34:06.783218",
    "2021-10-05 14:34:06.783221",
    "2021-10-05 14:34:06.783224",
    "2021-10-05 14:34:06.783226",
    "2021-10-05 14:34:06.783229",
    "2021-10-05 14:34:06.783232",
    "2021-10-05 14:34:06.783234",
    "2021-10-05 14:34:06.783237",
    "2021-10-05 14:34:06.783239",
    "2021-10-05 14:34:06.783242",
    "2021-10-05 14:34:06.783244",
    "2021-10-05 14:34:06.783247",
    "2021-10-05 14:34:06.783249",
    "2021-10-05 14:34:06.783252",
    "2021-10-05 14:34:06.783254",
    "2021-10-05 14:34:06.783257",
    "2021-10-05 14:34:06.783259",
    "2021-10-05 14:34:06.783262",
    "2021-10-05 14:34:06.783265",
    "2021-10-05 14:34:06.783267",
    "2021-10-05 14:34:06.783270",
    "2021-10-05 14:34:06.783272",
    "2021-10-05 14:34:06.783275",
    "2021-10-05 14:34:06.783277",
    "2021-10-05 14:34:06.783280",
    "2021-10-05 14:34:06.783282",
    "2021-10-05 14:34:06.783285",
    "2021-10-05 14:34:06.783288",
    "2021-10-05 14:34:06.783290",
    "2021-10-05 14:34:06.783293",
    "2021-10-05 14:34:06.783295",
    "2021-10-05 14:34:06.783297",
    "2021-10-05 14:34:06.783300",
    "2021-10-05 14:34:06.783302",
    "2021-10-05 14:34:06.783305",
    "2021-10-05 14:34:06.783308",
    "2021-10-05 14:34:06.783310",
    "2021-10-05 14:34:06.783313",
    "2021-10-05 14:34:06.783315",
    "2021-10-05 14:34:06.783318",
    "2021-10-05 14:34:06.783320",
    "2021-10-05 14:34:06.783323",
    "2021-10-05 14:34:06.783325",
    "2021-10-05 14:34:06.783328"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timeValues[$i]
}

# --- 2. Add a new "metadata" worksheet after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Hypertrophic cardiomyopathy_HCM"
$metaSheet.Range("C2").Value = 111
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.157"
$metaSheet.Range("D2").ClearFormats()
$metaSheet.Range("E2").Value = "2021-07-07T05:00:02.084745Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:34:06.779472"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/111/?format=json"

# Copy header styling (bold + border + centered) from the "data" sheet onto the
# new sheet's header row and the index column cell, reusing the same style.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("A1").Select()

Write-Output "done"
